$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 7 new data rows (233-239) below the existing data (which ends at 232).
# ---------------------------------------------------------------------------

# Copy the formatting (incl. the date number format on column A) from the
# last existing row down onto the new rows so the new cells pick up the same
# styles (e.g. style index used for date cells in column A).
$ws.Range("A232:T232").Copy() | Out-Null
$ws.Range("A233:T239").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row data: Date, Plant_Type, Plant_Size, Low, High, Rain, Growth, Quadrant,
#           Shade, UV, Humidity, Dew_Point, Pressure, Wind_Gust, Cloud_Cover,
#           Visibility, AQI, Pollen
$rows = @(
    @(233, 45820, "Flowering",    "Large",  66, 85, 0, 0,    "No", 2, "Bright",  8, 0.5, 64, 30.08, 10, 0.08, 9.9, 73, 46),
    @(234, 45820, "Nonflowering", "Medium", 66, 85, 0, 0,    "No", 3, "Bright",  8, 0.5, 64, 30.08, 10, 0.08, 9.9, 73, 46),
    @(235, 45820, "Nonflowering", "Small",  66, 85, 0, 0,    "No", 3, "Neutral", 8, 0.5, 64, 30.08, 10, 0.08, 9.9, 73, 46),
    @(236, 45820, "Nonflowering", "Medium", 66, 85, 0, 0.1,  "No", 3, "Bright",  8, 0.5, 64, 30.08, 10, 0.08, 9.9, 73, 46),
    @(237, 45820, "Nonflowering", "Medium", 66, 85, 0, 0.1,  "No", 3, "Bright",  8, 0.5, 64, 30.08, 10, 0.08, 9.9, 73, 46),
    @(238, 45820, "Nonflowering", "Large",  66, 85, 0, 0.25, "No", 4, "Bright",  8, 0.5, 64, 30.08, 10, 0.08, 9.9, 73, 46),
    @(239, 45820, "Tree",         "Medium", 66, 85, 0, 0.75, "No", 1, "Bright",  8, 0.5, 64, 30.08, 10, 0.08, 9.9, 73, 46)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]    # A Date
    $ws.Cells.Item($r, 2).Value = $row[2]    # B Plant_Type
    $ws.Cells.Item($r, 3).Value = $row[3]    # C Plant_Size
    $ws.Cells.Item($r, 4).Value = $row[4]    # D Low
    $ws.Cells.Item($r, 5).Value = $row[5]    # E High
    # F Temp_Diff formula is filled in below as a shared formula.
    $ws.Cells.Item($r, 7).Value = $row[6]    # G Rain
    $ws.Cells.Item($r, 8).Value = $row[7]    # H Growth
    $ws.Cells.Item($r, 9).Value = $row[8]    # I Quadrant
    $ws.Cells.Item($r, 10).Value = $row[9]   # J Shade
    $ws.Cells.Item($r, 11).Value = $row[10]  # K UV
    $ws.Cells.Item($r, 12).Value = $row[11]  # L Humidity
    $ws.Cells.Item($r, 13).Value = $row[12]  # M Dew_Point
    $ws.Cells.Item($r, 14).Value = $row[13]  # N Pressure
    $ws.Cells.Item($r, 15).Value = $row[14]  # O Wind_Gust
    $ws.Cells.Item($r, 16).Value = $row[15]  # P Cloud_Cover
    $ws.Cells.Item($r, 17).Value = $row[16]  # Q Visibility
    $ws.Cells.Item($r, 18).Value = $row[17]  # R AQI
    $ws.Cells.Item($r, 19).Value = $row[18]  # S Pollen
    $ws.Cells.Item($r, 20).Value = $row[19]  # T (last column)
}

# Extend the shared "Temp_Diff" formula (=ABS(D-E)) down through row 239.
$ws.Range("F233:F239").Formula = "=ABS(D233-E233)"

# ---------------------------------------------------------------------------
# Update the active selection to match the new view (cell U2 selected, sheet
# scrolled back to the top-left).
# ---------------------------------------------------------------------------
$ws.Range("U2").Select() | Out-Null
